$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.645.25'
$ws.Range("E2").Value = '  +1.25%  '

$ws.Range("D3").Value = '1.632.40'
$ws.Range("E3").Value = '  +1.39%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("E6").Value = '  +1.53%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +0.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.61%  '

$ws.Range("E10").Value = '  +2.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = '1.859.11'
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").Value = '1.638.10'
$ws.Range("E13").Value = '  +1.76%  '

$ws.Range("E14").Value = '  +1.53%  '

$ws.Range("E15").Value = '  +1.85%  '

$ws.Range("D16").Value = '26.633.49'
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.34%  '

$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").Value = '  +1.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.85%  '

$ws.Range("E20").Value = '  -0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.60%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.88%  '

$ws.Range("E24").Value = '  +2.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.04'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.46%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("E27").Value = '  -0.40%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.92%  '

$ws.Range("E30").Value = '  +4.94%  '

$ws.Range("E31").Value = '  -0.22%  '

$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.36'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.74%  '

$ws.Range("D36").Value = '1.167.94'
$ws.Range("E36").Value = '  +0.37%  '

$ws.Range("E37").Value = '  +1.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.810'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.38%  '

$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.504'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.51%  '

$ws.Range("E41").Value = '  -0.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.792'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").Value = '1.770.51'
$ws.Range("E44").Value = '  +1.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.33'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("E46").Value = '  +0.78%  '

$ws.Range("E47").Value = '  +1.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0511'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.83%  '

$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.409'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.43%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.98%  '

$ws.Range("E51").Value = '  -0.07%  '
